$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 80139019
$ws.Range("B2").Value = 88921
$ws.Range("E2").Value = 5741
$ws.Range("F2").Value = "Tjockfotad fingersvamp"
$ws.Range("G2").Value = "Ramaria flavescens"
$ws.Range("H2").Value = "(Schaeff.) R. H. Petersen"

$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"

$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("P2").Value = "Tveta friluftsgård, 300 m V om, Srm"
$ws.Range("Q2").Value = 648222.682956806
$ws.Range("R2").Value = 6560420.292955686
$ws.Range("S2").Value = 50

$ws.Range("Y2").Value = "'2019-09-27"
$ws.Range("AA2").Value = "'2019-09-27"

$ws.Range("AF2").ClearContents()

$ws.Range("AI2").Value = "barrskog"

$ws.Range("AW2").Value = "Hans Rydberg"
$ws.Range("AX2").Value = "Hans Rydberg"
